# Apply the data rotation/swap described by the diff.
# Rows 19, 20, 21 are cyclically rotated (19<-20, 20<-21, 21<-19)
# and rows 23, 24 are swapped (23<-24, 24<-23).
# Only columns A, B, E, F, G, H, Q, R, Z, AB change; everything else stays.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns touched by the edit, in order.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value2
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $vals[$col]
    }
}

# Capture original values before any writes (rows 19, 20, 21, 23, 24).
$row19 = Get-RowValues $ws 19 $cols
$row20 = Get-RowValues $ws 20 $cols
$row21 = Get-RowValues $ws 21 $cols
$row23 = Get-RowValues $ws 23 $cols
$row24 = Get-RowValues $ws 24 $cols

# Apply the rotation: new19 = old20, new20 = old21, new21 = old19.
Set-RowValues $ws 19 $cols $row20
Set-RowValues $ws 20 $cols $row21
Set-RowValues $ws 21 $cols $row19

# Apply the swap: new23 = old24, new24 = old23.
Set-RowValues $ws 23 $cols $row24
Set-RowValues $ws 24 $cols $row23
